$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 151 - this shifts the existing rows 151..162
# down to 152..163, matching the diff (each old row r's data now lives at r+1).
$ws.Rows("151").Insert()

# Populate the newly inserted row 151 with the new price-record data point.
$ws.Cells.Item(151, 1).Value = 5
$ws.Cells.Item(151, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(151, 3).Value = 'Maule'
$ws.Cells.Item(151, 4).Value = 44461
$ws.Cells.Item(151, 5).Value = 7
$ws.Cells.Item(151, 6).Value = 100112003
$ws.Cells.Item(151, 7).Value = 'Ajo'
$ws.Cells.Item(151, 8).Value = 'Chino'
$ws.Cells.Item(151, 9).Value = 'Primera'
$ws.Cells.Item(151, 10).Value = 200
$ws.Cells.Item(151, 11).Value = 15000
$ws.Cells.Item(151, 12).Value = 15000
$ws.Cells.Item(151, 13).Value = 15000
$ws.Cells.Item(151, 14).Value = '$/malla 10 kilos'
$ws.Cells.Item(151, 15).Value = 'China'
$ws.Cells.Item(151, 16).Value = 1500
$ws.Cells.Item(151, 17).Value = 10
$ws.Cells.Item(151, 18).Value = 'Hortaliza'

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(151, 4).NumberFormat = $ws.Cells.Item(152, 4).NumberFormat
